# Add files via upload
# The CPU sheet originally holds 4 "trio" blocks of FPS columns (1080p/1440p/4K)
# spanning B:T (B:D, E:G, H:J, K:M, N, O:Q, R:T). The author selected the last
# trio block H1:T16 (13 columns x 16 rows) and pasted it twice more to the
# right, extending the table from A1:T16 out to A1:AT16 (U1:AG16 = 1st paste,
# AH1:AT16 = 2nd paste).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPU")
$ws.Activate()

# Copy the source block (H1:T16) and paste it starting at U1.
$ws.Range("H1:T16").Copy()
$ws.Range("U1").PasteSpecial()

# Paste the very same block a second time, right after the first copy,
# starting at AH1 (U1:AG16 is 13 cols wide, so the next block starts at AH1).
$ws.Range("H1:T16").Copy()
$ws.Range("AH1").PasteSpecial()

# Clear the marching-ants clipboard marquee, like Excel does after pasting.
$excel.CutCopyMode = $false

# Mirror the final view state recorded in the workbook: the whole AH:AT
# column range ends up selected (this is what Excel leaves selected right
# after the second paste operation finishes).
$null = $ws.Range("AH1:AT1048576").Select()
